$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed by Excel as a number
# (e.g. "19.00" -> 19). Force them to stay plain text, matching the original
# inlineStr cells, then strip the temporary Text format so no extra cell style
# is left behind on the cell.
$textCells = @("D5", "D6", "D8", "D10", "D12", "D14", "D17", "D18", "D20", "D22", "D23", "D26", "D27", "D29", "D30", "D31", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D43", "D45", "D46", "D47", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D5").Value = "586.96"
$ws.Range("D6").Value = "189.36"
$ws.Range("D8").Value = "0.622"
$ws.Range("D10").Value = "0.215"
$ws.Range("D12").Value = "54.28"
$ws.Range("D14").Value = "9.42"
$ws.Range("D17").Value = "12.76"
$ws.Range("D18").Value = "19.00"
$ws.Range("D20").Value = "575.48"
$ws.Range("D22").Value = "0.996"
$ws.Range("D23").Value = "17.92"
$ws.Range("D26").Value = "94.36"
$ws.Range("D27").Value = "11.03"
$ws.Range("D29").Value = "9.34"
$ws.Range("D30").Value = "32.38"
$ws.Range("D31").Value = "7.08"
$ws.Range("D33").Value = "0.115"
$ws.Range("D34").Value = "3.82"
$ws.Range("D35").Value = "63.28"
$ws.Range("D36").Value = "3.22"
$ws.Range("D37").Value = "528.64"
$ws.Range("D38").Value = "0.407"
$ws.Range("D39").Value = "38.29"
$ws.Range("D43").Value = "3.53"
$ws.Range("D45").Value = "0.0457"
$ws.Range("D46").Value = "2.96"
$ws.Range("D47").Value = "3.47"
$ws.Range("D50").Value = "0.999"
$ws.Range("D51").Value = "1.43"

foreach ($addr in $textCells) {
    $ws.Range($addr).ClearFormats()
}

# Remaining cells: plain text already (non-numeric-looking), safe to assign directly.
$ws.Range("D2").Value = "70.681.88"
$ws.Range("E2").Value = "  +1.98%  "
$ws.Range("D3").Value = "3.565.79"
$ws.Range("E3").Value = "  +2.32%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("E5").Value = "  +1.18%  "
$ws.Range("E6").Value = "  +4.03%  "
$ws.Range("D7").Value = "3.561.57"
$ws.Range("E7").Value = "  +2.54%  "
$ws.Range("E8").Value = "  +2.20%  "
$ws.Range("E9").Value = "  -0.20%  "
$ws.Range("E10").Value = "  +10.40%  "
$ws.Range("E11").Value = "  +0.96%  "
$ws.Range("E12").Value = "  +1.32%  "
$ws.Range("E13").Value = "  +3.09%  "
$ws.Range("D15").Value = "4.134.64"
$ws.Range("E15").Value = "  +2.48%  "
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "70.700.88"
$ws.Range("E16").Value = "  +2.11%  "
$ws.Range("B17").Value = "Uniswap"
$ws.Range("C17").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("E17").Value = "  +4.27%  "
$ws.Range("B18").Value = "Chainlink"
$ws.Range("C18").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("E18").Value = "  -0.90%  "
$ws.Range("B19").Value = "WrappedEther"
$ws.Range("C19").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D19").Value = "3.534.29"
$ws.Range("E19").Value = "  +1.40%  "
$ws.Range("E20").Value = "  +7.80%  "
$ws.Range("E21").Value = "  +0.78%  "
$ws.Range("E22").Value = "  -0.82%  "
$ws.Range("E23").Value = "  -2.54%  "
$ws.Range("E24").Value = "  +3.11%  "
$ws.Range("E25").Value = "  +0.76%  "
$ws.Range("E26").Value = "  -1.03%  "
$ws.Range("E27").Value = "  +0.37%  "
$ws.Range("E28").Value = "  -0.85%  "
$ws.Range("E29").Value = "  +3.42%  "
$ws.Range("E30").Value = "  +2.21%  "
$ws.Range("E31").Value = "  -1.29%  "
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("E33").Value = "  +2.31%  "
$ws.Range("B34").Value = "dogwifhat"
$ws.Range("C34").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("E34").Value = "  +26.91%  "
$ws.Range("E35").Value = "  -0.43%  "
$ws.Range("E36").Value = "  +5.06%  "
$ws.Range("B37").Value = "Bittensor"
$ws.Range("C37").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("E37").Value = "  -0.07%  "
$ws.Range("B38").Value = "TheGraph"
$ws.Range("C38").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("E38").Value = "  +0.52%  "
$ws.Range("E39").Value = "  +1.24%  "
$ws.Range("D40").Value = "3.679.34"
$ws.Range("E40").Value = "  +10.30%  "
$ws.Range("D42").Value = "0.0₃0788"
$ws.Range("E42").Value = "  +4.76%  "
$ws.Range("E43").Value = "  +4.89%  "
$ws.Range("E44").Value = "  +3.64%  "
$ws.Range("E45").Value = "  +4.71%  "
$ws.Range("B46").Value = "ThetaToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("E46").Value = "  +0.27%  "
$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("E47").Value = "  -0.23%  "
$ws.Range("E48").Value = "  +3.43%  "
$ws.Range("E49").Value = "  +3.62%  "
$ws.Range("E50").Value = "  +0.01%  "
$ws.Range("E51").Value = "  +7.80%  "
